$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New ticker list for column A on Sheet1 (futures added: ES, NQ, XU front/next month; AA front+next)
$tickers = @(
    "Ident (BB_TCM)",
    "000333 C2 Equity",
    "600584 C1 Equity",
    "XUV1 Index",
    "HCTV1 Index",
    "ESZ1 Index",
    "1060 HK Equity",
    "XUU1 Index",
    "600519 CH Equity",
    "XUV1 Index",
    "11 HK Equity",
    "HCTV1 Index",
    "763 HK Equity",
    "5 HK Equity",
    "BABA US Equity",
    "KWEB US Equity",
    "LU US Equity",
    "HSCEI Index",
    "SHSZ300 Index",
    "XIN9I Index",
    "NDX Index",
    "HIU1 Index"
)

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 1
    $ws1.Range("A$row").Value = $tickers[$i]
}

# Update selections to match the saved view state
$null = $ws1.Range("E9").Select()
$null = $ws2.Range("D8").Select()

$null = $ws1.Activate()
